# Applies the "finish the fishery.estimates loop code for SKJ OBJ" edit.
# Strategy: use Range.InsertXML to rewrite the runs of a few existing
# paragraphs (splitting plain <w:t> runs into proofErr-wrapped runs, the
# way Word's spell/grammar checker would after a manual retype), mark the
# picture paragraph's run NoProofing (=> <w:rPr><w:noProof/></w:rPr>), and
# append four brand-new paragraphs at the end of the story.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) "lfmm <- read.lfmmdata.f(raw_data_dir,"LengthMM2000-2021.txt")"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(4)
$xml1 = '<w:p ' + $wNs + '>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>lfmm</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> &lt;- </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>read.lfmmdata</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>.f</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>(raw_data_dir,"LengthMM2000-2021.txt")</w:t></w:r>' +
  '</w:p>'
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "lfgrpd <- read.lengthfreq.f(raw_data_dir,"LengthFreq2000-2021.txt")"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(6)
$xml2 = '<w:p ' + $wNs + '>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>lfgrpd</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> &lt;- </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>read.lengthfreq</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>.f</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>(raw_data_dir,"LengthFreq2000-2021.txt")</w:t></w:r>' +
  '</w:p>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "What is month.substitution.mat?" (keeps the lastRenderedPageBreak)
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(8)
$xml3 = '<w:p ' + $wNs + '>' +
  '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">What is </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>month.substitution.mat</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>?</w:t></w:r>' +
  '</w:p>'
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) Picture paragraph: mark the run NoProofing so Word emits
#    <w:rPr><w:noProof/></w:rPr> ahead of <w:drawing>.
# ---------------------------------------------------------------------
$picPara = $d.Paragraphs($d.Paragraphs.Count)
$picPara.Range.NoProofing = $true

# ---------------------------------------------------------------------
# 5) Append the new trailing paragraphs after the picture paragraph:
#    - a page break
#    - "get.catch.estimates.f = function(...) {"
#    - an empty paragraph
#    - "why minsamps.in is set to 2?"
# ---------------------------------------------------------------------
$picPara.Range.InsertParagraphAfter()
$breakPara = $d.Paragraphs($d.Paragraphs.Count)
$breakPara.Range.InsertXML('<w:p ' + $wNs + '><w:r><w:br w:type="page"/></w:r></w:p>')

$breakPara.Range.InsertParagraphAfter()
$fnPara = $d.Paragraphs($d.Paragraphs.Count)
$xmlFn = '<w:p ' + $wNs + '>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>get.catch</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t>.estimates.f</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> = function(cae.in,caestrtflg.in,totunlds.in,lfgrpd.in,lfgrpd.stratflg.in,lfmm.in,my.year,minsamps.in,well.estimates,myarea.submat,growshrink.incrs.mat.touse,PS,Species) {</w:t></w:r>' +
  '</w:p>'
$fnPara.Range.InsertXML($xmlFn)

$fnPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.InsertXML('<w:p ' + $wNs + '/>')
$blankPara = $d.Paragraphs($d.Paragraphs.Count)

$blankPara.Range.InsertParagraphAfter()
$whyPara = $d.Paragraphs($d.Paragraphs.Count)
$xmlWhy = '<w:p ' + $wNs + '>' +
  '<w:r><w:t xml:space="preserve">why </w:t></w:r>' +
  '<w:r><w:t>minsamps.in</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> is set to 2</w:t></w:r>' +
  '<w:r><w:t>?</w:t></w:r>' +
  '</w:p>'
$whyPara.Range.InsertXML($xmlWhy)
